$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Target cluster labels:
# Row 4 was "Inflammatory-Mac" -> becomes "MuSCs"
# Row 5 was "MuSCs" -> becomes "Resolving-Mac"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "Resolving-Mac"

# Remove the old row 6 (previously "Resolving-Mac" target row) entirely
$ws.Rows.Item(6).Delete()

# Update recomputed TPM-based numeric values for rows 2-5

# Row 2 (Target cluster: ECs)
$ws.Range("G2").Value = 0.7999296666666668
$ws.Range("H2").Value = 2.399789
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.322136
$ws.Range("N2").Value = 0.966408
$ws.Range("O2").Value = 0.07882597898094613
$ws.Range("P2").Value = 0.07882597898094612
$ws.Range("Q2").Value = 0.2576861431013334
$ws.Range("R2").Value = 2.319175287912
$ws.Range("S2").Value = 0.07882597898094613
$ws.Range("T2").Value = 0.07882597898094612

# Row 3 (Target cluster: FAPs)
$ws.Range("G3").Value = 0.7999296666666668
$ws.Range("H3").Value = 2.399789
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.30984
$ws.Range("N3").Value = 9.92952
$ws.Range("O3").Value = 0.8099106534826741
$ws.Range("P3").Value = 0.8099106534826741
$ws.Range("Q3").Value = 2.64763920792
$ws.Range("R3").Value = 23.82875287128
$ws.Range("S3").Value = 0.8099106534826741
$ws.Range("T3").Value = 0.8099106534826741

# Row 4 (Target cluster: MuSCs)
$ws.Range("G4").Value = 0.7999296666666668
$ws.Range("H4").Value = 2.399789
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.438826
$ws.Range("N4").Value = 1.316478
$ws.Range("O4").Value = 0.1073797683347799
$ws.Range("P4").Value = 0.1073797683347799
$ws.Range("Q4").Value = 0.3510299359046667
$ws.Range("R4").Value = 3.159269423142
$ws.Range("S4").Value = 0.1073797683347799
$ws.Range("T4").Value = 0.1073797683347799

# Row 5 (Target cluster: Resolving-Mac)
$ws.Range("G5").Value = 0.7999296666666668
$ws.Range("H5").Value = 2.399789
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.015871
$ws.Range("N5").Value = 0.047613
$ws.Range("O5").Value = 0.003883599201599933
$ws.Range("P5").Value = 0.003883599201599933
$ws.Range("Q5").Value = 0.01269568373966667
$ws.Range("R5").Value = 0.114261153657
$ws.Range("S5").Value = 0.003883599201599933
$ws.Range("T5").Value = 0.003883599201599933
